$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 114, shifting existing rows 114-159 down to 115-160.
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with the new weekly price record.
$ws.Range("A114").Value = 10
$ws.Range("B114").Value = "Vega Modelo de Temuco"
$ws.Range("C114").Value = "La Araucanía"
$ws.Range("D114").Value = 44795
$ws.Range("E114").Value = 9
$ws.Range("F114").Value = "Fruta"
$ws.Range("G114").Value = 100104
$ws.Range("H114").Value = "Frutos de pepita"
$ws.Range("I114").Value = 100104001
$ws.Range("J114").Value = "Granada"
$ws.Range("K114").Value = "Wonderfull"
$ws.Range("L114").Value = "Primera"
$ws.Range("M114").Value = 55
$ws.Range("N114").Value = 14000
$ws.Range("O114").Value = 15000
$ws.Range("P114").Value = 14455
$ws.Range("Q114").Value = "$/bandeja 10 kilos granel"
$ws.Range("R114").Value = "Provincia de Limarí"
$ws.Range("S114").Value = 1446
$ws.Range("T114").Value = 10
